# Workbook edit: add JSON `storeKeys(json,jsonpath,var)` keyword, and remove the
# stray single-entry "text" category (both its row in the `target` list and its
# data column) from the hidden "#system" sheet that backs the macro editor's
# autocomplete lists.
#
# NOTE: this engine's Range.Insert/Delete with xlShiftDown/xlShiftUp shifts the
# *entire row* (all columns), not just the column of the target range, so we
# shift the handful of affected cells manually instead of relying on Insert/
# Delete for the vertical moves. (Whole-column Delete works fine and is used
# for the column removal below.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) JSON keywords (column M): insert "storeKeys(json,jsonpath,var)" in its
#    alphabetically-correct slot (between storeCount and storeValue), pushing
#    storeValue/storeValues down by one row. Walk bottom-up so we don't clobber
#    a value before it's copied down.
# ---------------------------------------------------------------------------
$ws.Range("M18").Value = $ws.Range("M17").Value2
$ws.Range("M17").Value = $ws.Range("M16").Value2
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 2) Category list (column A): drop the "text" entry (row 25) entirely,
#    shifting the remaining categories (web, webalert, webcookie, ws,
#    ws.async, xml) up by one row, then clear the now-vacated last row.
# ---------------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $ws.Range("A$r").Value = $ws.Range("A$($r+1)").Value2
}
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------------
# 3) Data column for "text" (column Y) is now orphaned/empty -- delete the
#    whole column so web/webalert/webcookie/ws/ws.async/xml each shift one
#    column to the left. (Whole-column delete is correctly scoped.)
# ---------------------------------------------------------------------------
$ws.Columns("Y").Delete()

# ---------------------------------------------------------------------------
# 4) Fix up the workbook-level defined names that point at the ranges we just
#    resized/moved (this engine does not auto-repair names on structural
#    edits the way interactive Excel does).
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    switch ($n.Name) {
        "json"      { $n.RefersTo = "='#system'!`$M`$2:`$M`$18" }
        "target"    { $n.RefersTo = "='#system'!`$A`$2:`$A`$30" }
        "web"       { $n.RefersTo = "='#system'!`$Y`$2:`$Y`$129" }
        "webalert"  { $n.RefersTo = "='#system'!`$Z`$2:`$Z`$8" }
        "webcookie" { $n.RefersTo = "='#system'!`$AA`$2:`$AA`$8" }
        "ws"        { $n.RefersTo = "='#system'!`$AB`$2:`$AB`$17" }
        "ws.async"  { $n.RefersTo = "='#system'!`$AC`$2:`$AC`$8" }
        "xml"       { $n.RefersTo = "='#system'!`$AD`$2:`$AD`$27" }
    }
}
